# Update gh-pages to output generated at a3196b5
#
# The "苏州·世纪幻想动漫游戏展" entry (old row 2 / record #1, dated
# 2024-01-21) was removed from both the "展览" and "全部类型" sheets.
# That shifts every subsequent record up by one row (so the sheet shrinks
# from A1:J24 to A1:J23), the "序号" numbers in column A are renumbered
# 1..22, and the "想去人数" (column F) counters were refreshed with newer
# live figures for all of the remaining 22 records.

# New column-F ("想去人数") values, in order, for the 22 records that
# remain after the first record is removed.
$newWantCounts = @(270, 274, 10631, 9336, 579, 678, 92, 26, 15, 9446, 2419, 25, 62, 358, 10675, 10703, 4, 8, 3, 4, 8, 4)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # 1. Remove the first data row (row 2) -- this shifts all following
    #    rows up by one and naturally carries along their formatting.
    $ws.Rows(2).Delete()

    # 2. Renumber the "序号" column (A) sequentially: header stays 0,
    #    data rows become 1..22.
    for ($i = 1; $i -le 22; $i++) {
        $ws.Cells.Item($i + 1, 1).Value = $i
    }

    # 3. Refresh the "想去人数" column (F) with the updated counts.
    for ($i = 0; $i -lt $newWantCounts.Length; $i++) {
        $ws.Cells.Item($i + 2, 6).Value = $newWantCounts[$i]
    }
}
